$wb = $excel.ActiveWorkbook

# Rename the existing sheet and add the new one right after it.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TwoPqAndOnePvNodeDifferentOrder"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TwoPqNodes"

# --- block 1 (row 1) ---
$ws2.Range("A1").Value = "admittances"
$ws2.Range("D1").Value = "magnitudes"
$ws2.Range("G1").Value = "angles"

$ws2.Range("A2").Formula = "=COMPLEX(1, 2)"
$ws2.Range("B2").Formula = "=COMPLEX(3, 4)"
$ws2.Range("D2").Formula = "=IMABS(A2)"
$ws2.Range("E2").Formula = "=IMABS(B2)"
$ws2.Range("G2").Formula = "=IMARGUMENT(A2)"
$ws2.Range("H2").Formula = "=IMARGUMENT(B2)"

$ws2.Range("A3").Formula = "=COMPLEX(5, 6)"
$ws2.Range("B3").Formula = "=COMPLEX(7,8)"
$ws2.Range("D3").Formula = "=IMABS(A3)"
$ws2.Range("E3").Formula = "=IMABS(B3)"
$ws2.Range("G3").Formula = "=IMARGUMENT(A3)"
$ws2.Range("H3").Formula = "=IMARGUMENT(B3)"

# --- voltages block ---
$ws2.Range("A5").Value = "voltages"
$ws2.Range("C5").Value = "magnitudes"
$ws2.Range("E5").Value = "angles"

$ws2.Range("A6").Formula = "=COMPLEX(9,10)"
$ws2.Range("C6").Formula = "=IMABS(A6)"
$ws2.Range("E6").Formula = "=IMARGUMENT(A6)"

$ws2.Range("A7").Formula = "=COMPLEX(11,12)"
$ws2.Range("C7").Formula = "=IMABS(A7)"
$ws2.Range("E7").Formula = "=IMARGUMENT(A7)"

# --- block 3 (row 9) ---
$ws2.Range("A9").Value = "currents"
$ws2.Range("C9").Value = "magnitudes"
$ws2.Range("E9").Value = "angles"

$ws2.Range("A10").Formula = "=COMPLEX(13,14)"
$ws2.Range("C10").Formula = "=IMABS(A10)"
$ws2.Range("E10").Formula = "=IMARGUMENT(A10)"

$ws2.Range("A11").Formula = "=COMPLEX(15,16)"
$ws2.Range("C11").Formula = "=IMABS(A11)"
$ws2.Range("E11").Formula = "=IMARGUMENT(A11)"

# --- real power by angle ---
$ws2.Range("A13").Value = "real power by angle"

$ws2.Range("A14").Formula = "=-C6*E2*C7*SIN(E6-H2-E7)+C10*C6*SIN(E6-E10)"
$ws2.Range("B14").Formula = "=C6*E2*C7*SIN(E6-H2-E7)"

$ws2.Range("A15").Formula = "=C6*D3*C7*SIN(E7-G3-E6)"
$ws2.Range("B15").Formula = "=-C7*D3*C6*SIN(E7-G3-E6)+C11*C7*SIN(E7-E11)"

# --- imaginary power by amplitude ---
$ws2.Range("A17").Value = "imaginary power by amplitude"

$ws2.Range("A18").Formula = "=E2*C7*SIN(E6-H2-E7)-2*D2*C6*SIN(G2)-C10*SIN(E6-E10)"
$ws2.Range("B18").Formula = "=C6*E2*SIN(E6-H2-E7)"

$ws2.Range("A19").Formula = "=C7*D3*SIN(E7-G3-E6)"
$ws2.Range("B19").Formula = "=D3*C6*SIN(E7-G3-E6)-2*E3*C7*SIN(H3)-C11*SIN(E7-E11)"

# Match the recorded selections: sheet1 keeps A2:B3 selected (not active),
# sheet2 ends up active with J19 selected.
[void]$ws1.Range("A2:B3").Select()
[void]$ws2.Range("J19").Select()
